# Update "想去人数" (F column) counts in the scraped convention-info workbook.
# Mirrors the regenerated output data (gh-pages build) described by the commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2695
$ws1.Range("F27").Value = 1377
$ws1.Range("F28").Value = 284
$ws1.Range("F32").Value = 357
$ws1.Range("F39").Value = 1426

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 81

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2695
$ws4.Range("F21").Value = 81
$ws4.Range("F27").Value = 1377
$ws4.Range("F28").Value = 284
$ws4.Range("F33").Value = 357
$ws4.Range("F43").Value = 1426
